$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "כאשר המערכת ריקה ולוחצים על שאילתה מקבלים Exception"
$ws.Range("B11").Value = "להוסיף יוזר ADMIN שלא נמצא על DB וניתן לעלות איתו תמיד"

$ws.Range("B12").Select()
